$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" contain the same event rows and both need updating.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 207
    $ws.Range("F5").Value = 55
}
